$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.098.55"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.840.75"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.37"
$ws.Range("E5").Value = "  -1.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6851"
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3020"
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07435"
$ws.Range("E9").Value = "  -3.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.09"
$ws.Range("E10").Value = "  -2.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07665"
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("D12").Value = "1.836.63"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.055"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6829"
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.42"
$ws.Range("E15").Value = "  -6.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.168"
$ws.Range("E16").Value = "  -7.07%  "
$ws.Range("D17").Value = "29.103.86"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008142"
$ws.Range("E18").Value = "  -2.26%  "
$ws.Range("D19").Value = "2.079.79"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.55"
$ws.Range("E20").Value = "  -6.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.53"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.382"
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.19"
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1450"
$ws.Range("E26").Value = "  -4.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.767"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.09"
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.513"
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.265"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.134"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05233"
$ws.Range("E33").Value = "  +2.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7579"
$ws.Range("E34").Value = "  -4.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.849"
$ws.Range("E35").Value = "  -3.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.133"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "1.306.35"
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01837"
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9322"
$ws.Range("E41").Value = "  -2.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.814"
$ws.Range("E42").Value = "  -4.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.81"
$ws.Range("E43").Value = "  -2.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.983.47"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000123"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5195"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.74"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.516"
$ws.Range("E49").Value = "  -2.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.770"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05942"
$ws.Range("E51").Value = "  +0.76%  "
